$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 7565.5
$ws.Range("I29").Value = 87.333336
$ws.Range("K29").Value = 262.000008
$ws.Range("M29").Value = 18.99999200000002
$ws.Range("H41").Value = 1344.8788
$ws.Range("I41").Value = 1271.2084
$ws.Range("K41").Value = 1271.2084
$ws.Range("M41").Value = -831.2084
$ws.Range("H55").Value = 137.15384
$ws.Range("I55").Value = 138.3
$ws.Range("J55").Value = 133.33333
$ws.Range("K55").Value = 138.3
$ws.Range("L55").Value = 133.33333
$ws.Range("M55").Value = 75.69999999999999
$ws.Range("N55").Value = -561.3333299999999
$ws.Range("I62").Value = 8505.471
$ws.Range("K62").Value = 8505.471
$ws.Range("M62").Value = -7881.471
$ws.Range("I65").Value = 8505.471
$ws.Range("K65").Value = 42527.355
$ws.Range("M65").Value = -39407.355
$ws.Range("H69").Value = 7583.5
$ws.Range("J69").Value = 6853.75
$ws.Range("L69").Value = 20561.25
$ws.Range("N69").Value = -22309.25
$ws.Range("H72").Value = 7583.5
$ws.Range("J72").Value = 6853.75
$ws.Range("L72").Value = 61683.75
$ws.Range("N72").Value = -70419.75
$ws.Range("H86").Value = 92848.63
$ws.Range("I86").Value = 168842.17
$ws.Range("J86").Value = 1656.4
$ws.Range("K86").Value = 168842.17
$ws.Range("L86").Value = 1656.4
$ws.Range("M86").Value = -167719.17
$ws.Range("N86").Value = -3902.4
$ws.Range("H89").Value = 92848.63
$ws.Range("I89").Value = 168842.17
$ws.Range("J89").Value = 1656.4
$ws.Range("K89").Value = 844210.8500000001
$ws.Range("L89").Value = 8282
$ws.Range("M89").Value = -838594.8500000001
$ws.Range("N89").Value = -19514
$ws.Range("H98").Value = 266624.12
$ws.Range("I98").Value = 1783.25
$ws.Range("K98").Value = 1783.25
$ws.Range("M98").Value = -285.25
$ws.Range("H100").Value = 1668.75
$ws.Range("I100").Value = 1208.9
$ws.Range("K100").Value = 1208.9
$ws.Range("M100").Value = -667.9000000000001
$ws.Range("H111").Value = 13049.23
$ws.Range("I111").Value = 16520.857
$ws.Range("K111").Value = 49562.571
$ws.Range("M111").Value = -46495.571
$ws.Range("H112").Value = 2946899.8
$ws.Range("J112").Value = 2946899.8
$ws.Range("L112").Value = 8840699.399999999
$ws.Range("N112").Value = -8842915.399999999
$ws.Range("H113").Value = 5440.1113
$ws.Range("I113").Value = 7538.8335
$ws.Range("J113").Value = 1242.6666
$ws.Range("K113").Value = 7538.8335
$ws.Range("L113").Value = 1242.6666
$ws.Range("M113").Value = -4284.8335
$ws.Range("N113").Value = -7750.6666
$ws.Range("H116").Value = 8691.103999999999
$ws.Range("I116").Value = 12030.6875
$ws.Range("J116").Value = 4580.846
$ws.Range("K116").Value = 12030.6875
$ws.Range("L116").Value = 4580.846
$ws.Range("M116").Value = -8588.6875
$ws.Range("N116").Value = -11464.846
$ws.Range("H122").Value = 266624.12
$ws.Range("I122").Value = 1783.25
$ws.Range("K122").Value = 5349.75
$ws.Range("M122").Value = -2899.75
$ws.Range("H138").Value = 2479.6304
$ws.Range("I138").Value = 1751.7222
$ws.Range("J138").Value = 2947.5715
$ws.Range("K138").Value = 5255.1666
$ws.Range("L138").Value = 8842.7145
$ws.Range("M138").Value = -115.1665999999996
$ws.Range("N138").Value = -19122.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1555.24
$ws.Range("I45").Value = 917.4286
$ws.Range("K45").Value = 917.4286
$ws.Range("M45").Value = -540.4286
$ws.Range("H61").Value = 7664.7
$ws.Range("J61").Value = 8441.333000000001
$ws.Range("L61").Value = 8441.333000000001
$ws.Range("N61").Value = -8865.333000000001
$ws.Range("H102").Value = 2790.7273
$ws.Range("I102").Value = 2343.0667
$ws.Range("J102").Value = 3750
$ws.Range("K102").Value = 2343.0667
$ws.Range("L102").Value = 3750
$ws.Range("M102").Value = -721.0666999999999
$ws.Range("N102").Value = -6994
$ws.Range("H110").Value = 6329.154
$ws.Range("I110").Value = 4722.5713
$ws.Range("K110").Value = 4722.5713
$ws.Range("M110").Value = -2677.5713
$ws.Range("H122").Value = 8659.083000000001
$ws.Range("I122").Value = 6153.294
$ws.Range("J122").Value = 14744.571
$ws.Range("K122").Value = 18459.882
$ws.Range("L122").Value = 44233.713
$ws.Range("M122").Value = -16009.882
$ws.Range("N122").Value = -49133.713
$ws.Range("H132").Value = 4684.921
$ws.Range("I132").Value = 2944.7058
$ws.Range("K132").Value = 8834.117400000001
$ws.Range("M132").Value = -6304.117400000001
$ws.Range("H136").Value = 7664.7
$ws.Range("J136").Value = 8441.333000000001
$ws.Range("L136").Value = 25323.999
$ws.Range("N136").Value = -30423.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3076.697
$ws.Range("J20").Value = 5019.8887
$ws.Range("L20").Value = 5019.8887
$ws.Range("N20").Value = -5513.8887
$ws.Range("H86").Value = 1822.0526
$ws.Range("I86").Value = 1658.6428
$ws.Range("J86").Value = 2279.6
$ws.Range("K86").Value = 1658.6428
$ws.Range("L86").Value = 2279.6
$ws.Range("M86").Value = -535.6428000000001
$ws.Range("N86").Value = -4525.6
$ws.Range("H89").Value = 1822.0526
$ws.Range("I89").Value = 1658.6428
$ws.Range("J89").Value = 2279.6
$ws.Range("K89").Value = 8293.214
$ws.Range("L89").Value = 11398
$ws.Range("M89").Value = -2677.214
$ws.Range("N89").Value = -22630
$ws.Range("H105").Value = 1624.1666
$ws.Range("I105").Value = 1567.6842
$ws.Range("J105").Value = 1838.8
$ws.Range("K105").Value = 1567.6842
$ws.Range("L105").Value = 1838.8
$ws.Range("M105").Value = 179.3158000000001
$ws.Range("N105").Value = -5332.8
$ws.Range("H107").Value = 2013.55
$ws.Range("I107").Value = 1971.129
$ws.Range("K107").Value = 1971.129
$ws.Range("M107").Value = -51.12899999999991
$ws.Range("H134").Value = 1669.697
$ws.Range("I134").Value = 1466.0377
$ws.Range("K134").Value = 4398.1131
$ws.Range("M134").Value = -1863.1131
$ws.Range("H140").Value = 69000
$ws.Range("J140").Value = 69000
$ws.Range("L140").Value = 69000
$ws.Range("N140").Value = -79360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1836
$ws.Range("I16").Value = 1610.8334
$ws.Range("K16").Value = 1610.8334
$ws.Range("M16").Value = -1323.8334
$ws.Range("H22").Value = 837.6667
$ws.Range("I22").Value = 520
$ws.Range("J22").Value = 1473
$ws.Range("K22").Value = 520
$ws.Range("L22").Value = 1473
$ws.Range("M22").Value = -170
$ws.Range("N22").Value = -2173
$ws.Range("H43").Value = 25928.8
$ws.Range("J43").Value = 25928.8
$ws.Range("L43").Value = 25928.8
$ws.Range("N43").Value = -26296.8
$ws.Range("H62").Value = 2861.8
$ws.Range("J62").Value = 2873.75
$ws.Range("L62").Value = 2873.75
$ws.Range("N62").Value = -4121.75
$ws.Range("H65").Value = 2861.8
$ws.Range("J65").Value = 2873.75
$ws.Range("L65").Value = 14368.75
$ws.Range("N65").Value = -20608.75
$ws.Range("H101").Value = 25928.8
$ws.Range("J101").Value = 25928.8
$ws.Range("L101").Value = 25928.8
$ws.Range("N101").Value = -32418.8
$ws.Range("H113").Value = 1836
$ws.Range("I113").Value = 1610.8334
$ws.Range("K113").Value = 1610.8334
$ws.Range("M113").Value = 559.1666
$ws.Range("H122").Value = 1813.2778
$ws.Range("I122").Value = 1418.3334
$ws.Range("J122").Value = 3788
$ws.Range("K122").Value = 4255.0002
$ws.Range("L122").Value = 11364
$ws.Range("M122").Value = -1805.0002
$ws.Range("N122").Value = -16264
$ws.Range("H132").Value = 3012.9722
$ws.Range("I132").Value = 2668.4517
$ws.Range("K132").Value = 8005.355100000001
$ws.Range("M132").Value = -5475.355100000001
$ws.Range("H134").Value = 5527.231
$ws.Range("I134").Value = 5995.391
$ws.Range("K134").Value = 17986.173
$ws.Range("M134").Value = -15451.173

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 659.86957
$ws.Range("I12").Value = 969
$ws.Range("K12").Value = 2907
$ws.Range("M12").Value = -2734
$ws.Range("H38").Value = 135.81818
$ws.Range("I38").Value = 144.5
$ws.Range("K38").Value = 433.5
$ws.Range("M38").Value = -86.5
$ws.Range("H113").Value = 1481.7273
$ws.Range("I113").Value = 1087.5
$ws.Range("J113").Value = 1569.3334
$ws.Range("K113").Value = 3262.5
$ws.Range("L113").Value = 4708.0002
$ws.Range("M113").Value = -1092.5
$ws.Range("N113").Value = -9048.0002
$ws.Range("H131").Value = 7354510
$ws.Range("J131").Value = 1772.6909
$ws.Range("L131").Value = 5318.072700000001
$ws.Range("N131").Value = -15398.0727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12499.667
$ws.Range("I70").Value = 11199.6
$ws.Range("K70").Value = 11199.6
$ws.Range("M70").Value = -10929.6
$ws.Range("H73").Value = 12499.667
$ws.Range("I73").Value = 11199.6
$ws.Range("K73").Value = 11199.6
$ws.Range("M73").Value = -10263.6
$ws.Range("H80").Value = 4184.8335
$ws.Range("I80").Value = 4011.625
$ws.Range("J80").Value = 4323.4
$ws.Range("K80").Value = 4011.625
$ws.Range("L80").Value = 4323.4
$ws.Range("M80").Value = -3013.625
$ws.Range("N80").Value = -6319.4
$ws.Range("H83").Value = 4184.8335
$ws.Range("I83").Value = 4011.625
$ws.Range("J83").Value = 4323.4
$ws.Range("K83").Value = 20058.125
$ws.Range("L83").Value = 21617
$ws.Range("M83").Value = -15066.125
$ws.Range("N83").Value = -31601
$ws.Range("H102").Value = 1303.4894
$ws.Range("I102").Value = 1362.5581
$ws.Range("K102").Value = 1362.5581
$ws.Range("M102").Value = 259.4419
$ws.Range("H113").Value = 9055.166999999999
$ws.Range("I113").Value = 9874.5625
$ws.Range("K113").Value = 9874.5625
$ws.Range("M113").Value = -7704.5625
$ws.Range("H122").Value = 1694.2142
$ws.Range("I122").Value = 1397.5
$ws.Range("J122").Value = 3474.5
$ws.Range("K122").Value = 4192.5
$ws.Range("L122").Value = 10423.5
$ws.Range("M122").Value = -1742.5
$ws.Range("N122").Value = -15323.5
$ws.Range("H132").Value = 28042.426
$ws.Range("I132").Value = 27963.078
$ws.Range("K132").Value = 83889.234
$ws.Range("M132").Value = -81359.234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2409.7778
$ws.Range("I16").Value = 1185.1364
$ws.Range("J16").Value = 7798.2
$ws.Range("K16").Value = 1185.1364
$ws.Range("L16").Value = 7798.2
$ws.Range("M16").Value = -1015.1364
$ws.Range("N16").Value = -8138.2
$ws.Range("H22").Value = 2122.2666
$ws.Range("I22").Value = 1258.75
$ws.Range("J22").Value = 3109.1428
$ws.Range("K22").Value = 1258.75
$ws.Range("L22").Value = 3109.1428
$ws.Range("M22").Value = -963.75
$ws.Range("N22").Value = -3699.1428
$ws.Range("H27").Value = 2122.2666
$ws.Range("I27").Value = 1258.75
$ws.Range("J27").Value = 3109.1428
$ws.Range("K27").Value = 1258.75
$ws.Range("L27").Value = 3109.1428
$ws.Range("M27").Value = -1151.75
$ws.Range("N27").Value = -3323.1428
$ws.Range("H42").Value = 36998
$ws.Range("I42").Value = 36998
$ws.Range("K42").Value = 36998
$ws.Range("M42").Value = -36435
$ws.Range("H46").Value = 2735.318
$ws.Range("I46").Value = 1222.8823
$ws.Range("J46").Value = 7877.6
$ws.Range("K46").Value = 1222.8823
$ws.Range("L46").Value = 7877.6
$ws.Range("M46").Value = -1034.8823
$ws.Range("N46").Value = -8253.6
$ws.Range("H49").Value = 36998
$ws.Range("I49").Value = 36998
$ws.Range("K49").Value = 36998
$ws.Range("M49").Value = -36851
$ws.Range("H68").Value = 2666.0303
$ws.Range("I68").Value = 2671.25
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 2671.25
$ws.Range("L68").Value = 2499
$ws.Range("M68").Value = -1922.25
$ws.Range("N68").Value = -3997
$ws.Range("H71").Value = 2666.0303
$ws.Range("I71").Value = 2671.25
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 13356.25
$ws.Range("L71").Value = 12495
$ws.Range("M71").Value = -9612.25
$ws.Range("N71").Value = -19983
$ws.Range("H106").Value = 23702
$ws.Range("J106").Value = 23702
$ws.Range("L106").Value = 23702
$ws.Range("N106").Value = -26226
$ws.Range("H123").Value = 20000
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -15100
$ws.Range("N123").ClearContents()
$ws.Range("H132").Value = 8372.277
$ws.Range("I132").Value = 1940.875
$ws.Range("J132").Value = 10209.821
$ws.Range("K132").Value = 5822.625
$ws.Range("L132").Value = 30629.463
$ws.Range("M132").Value = -3292.625
$ws.Range("N132").Value = -35689.463
$ws.Range("H136").Value = 2167.25
$ws.Range("I136").Value = 1750.9474
$ws.Range("J136").Value = 3749.2
$ws.Range("K136").Value = 5252.8422
$ws.Range("L136").Value = 11247.6
$ws.Range("M136").Value = -2702.8422
$ws.Range("N136").Value = -16347.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 39166.5
$ws.Range("J12").Value = 45000
$ws.Range("L12").Value = 45000
$ws.Range("N12").Value = -45284
$ws.Range("H51").Value = 17999.5
$ws.Range("J51").Value = 21999
$ws.Range("L51").Value = 21999
$ws.Range("N51").Value = -23019
$ws.Range("H52").Value = 30213.25
$ws.Range("I52").Value = 24499.5
$ws.Range("J52").Value = 35927
$ws.Range("K52").Value = 24499.5
$ws.Range("L52").Value = 35927
$ws.Range("M52").Value = -24273.5
$ws.Range("N52").Value = -36379
$ws.Range("H62").Value = 8518.056
$ws.Range("I62").Value = 5596.6665
$ws.Range("K62").Value = 5596.6665
$ws.Range("M62").Value = -4972.6665
$ws.Range("H65").Value = 8518.056
$ws.Range("I65").Value = 5596.6665
$ws.Range("K65").Value = 27983.3325
$ws.Range("M65").Value = -24863.3325
$ws.Range("H74").Value = 29996.334
$ws.Range("J74").Value = 29996.334
$ws.Range("L74").Value = 29996.334
$ws.Range("N74").Value = -31868.334
$ws.Range("H77").Value = 29996.334
$ws.Range("J77").Value = 29996.334
$ws.Range("L77").Value = 89989.00199999999
$ws.Range("N77").Value = -99349.00199999999
$ws.Range("H107").Value = 433.6154
$ws.Range("I107").Value = 426.66666
$ws.Range("K107").Value = 1279.99998
$ws.Range("M107").Value = 640.0000199999999
$ws.Range("H112").Value = 32947.2
$ws.Range("J112").Value = 35571.5
$ws.Range("L112").Value = 35571.5
$ws.Range("N112").Value = -38525.5
$ws.Range("H126").Value = 5447.2856
$ws.Range("I126").Value = 2499
$ws.Range("J126").Value = 6626.6
$ws.Range("K126").Value = 7497
$ws.Range("L126").Value = 19879.8
$ws.Range("M126").Value = -5027
$ws.Range("N126").Value = -24819.8
$ws.Range("H132").Value = 1974.7317
$ws.Range("I132").Value = 1054.5
$ws.Range("K132").Value = 3163.5
$ws.Range("M132").Value = -633.5
